$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (GitHub Actions scheduled refresh)

$ws.Cells.Item(2, 4).Value = '34.914.00'
$ws.Cells.Item(2, 5).Value = '  +1.44%  '
$ws.Cells.Item(3, 4).Value = '1.818.42'
$ws.Cells.Item(3, 5).Value = '  +1.10%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '231.01'
$ws.Cells.Item(5, 5).Value = '  +3.04%  '
$ws.Cells.Item(6, 5).Value = '  +1.97%  '
$ws.Cells.Item(7, 5).Value = '  +0.38%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '41.80'
$ws.Cells.Item(8, 5).Value = '  +1.06%  '
$ws.Cells.Item(9, 5).Value = '  +6.84%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0683'
$ws.Cells.Item(10, 5).Value = '  +2.71%  '
$ws.Cells.Item(11, 5).Value = '  +0.73%  '
$ws.Cells.Item(12, 4).Value = '2.085.16'
$ws.Cells.Item(12, 5).Value = '  +1.28%  '
$ws.Cells.Item(13, 4).Value = '1.838.00'
$ws.Cells.Item(13, 5).Value = '  +2.24%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '11.04'
$ws.Cells.Item(14, 5).Value = '  +2.11%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.658'
$ws.Cells.Item(15, 5).Value = '  +5.59%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '4.64'
$ws.Cells.Item(16, 5).Value = '  +6.29%  '
$ws.Cells.Item(17, 4).Value = '34.949.02'
$ws.Cells.Item(17, 5).Value = '  +1.50%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '69.39'
$ws.Cells.Item(18, 5).Value = '  +3.31%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0785'
$ws.Cells.Item(19, 5).Value = '  +2.75%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '237.99'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '11.74'
$ws.Cells.Item(21, 5).Value = '  +6.19%  '
$ws.Cells.Item(22, 5).Value = '  +0.37%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '4.56'
$ws.Cells.Item(23, 5).Value = '  +11.99%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.25'
$ws.Cells.Item(24, 5).Value = '  +4.59%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '172.07'
$ws.Cells.Item(25, 5).Value = '  +0.40%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '7.74'
$ws.Cells.Item(26, 5).Value = '  +1.73%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '17.37'
$ws.Cells.Item(28, 5).Value = '  -0.03%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.56'
$ws.Cells.Item(29, 5).Value = '  +27.07%  '
$ws.Cells.Item(30, 5).Value = '  +0.42%  '
$ws.Cells.Item(31, 4).Value = '3.348.75'
$ws.Cells.Item(31, 5).Value = '  +37.83%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.0548'
$ws.Cells.Item(32, 5).Value = '  +7.48%  '
$ws.Cells.Item(33, 5).Value = '  +3.01%  '
$ws.Cells.Item(34, 5).Value = '  +3.97%  '
$ws.Cells.Item(35, 5).Value = '  +1.09%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '92.17'
$ws.Cells.Item(36, 5).Value = '  +8.57%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.672'
$ws.Cells.Item(37, 5).Value = '  +4.55%  '
$ws.Cells.Item(38, 5).Value = '  +5.13%  '
$ws.Cells.Item(39, 4).Value = '1.313.28'
$ws.Cells.Item(39, 5).Value = '  -0.05%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0191'
$ws.Cells.Item(40, 5).Value = '  +2.44%  '
$ws.Cells.Item(41, 5).Value = '  +3.27%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.981'
$ws.Cells.Item(42, 5).Value = '  +5.01%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '14.67'
$ws.Cells.Item(43, 5).Value = '  +0.05%  '
$ws.Cells.Item(44, 2).Value = 'HuobiToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.45'
$ws.Cells.Item(44, 5).Value = '  +0.62%  '
$ws.Cells.Item(45, 2).Value = 'RenderToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.28'
$ws.Cells.Item(45, 5).Value = '  -1.52%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.77'
$ws.Cells.Item(46, 5).Value = '  -0.63%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '6.16'
$ws.Cells.Item(47, 5).Value = '  +6.08%  '
$ws.Cells.Item(48, 5).Value = '  -1.75%  '
$ws.Cells.Item(49, 4).Value = '1.997.78'
$ws.Cells.Item(49, 5).Value = '  +1.96%  '
$ws.Cells.Item(50, 5).Value = '  +0.40%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '100.04'
$ws.Cells.Item(51, 5).Value = '  -0.37%  '
